$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gene Cluster Instance - Report")

# 1. Update the base path used for all GeneCluster Genbank file references:
#    .../showcase_examples/Aflavus_aflatoxin_and_leporinB/leporinB_redo/...
#    -> .../showcase_examples_redo/Aflavus_aflatoxin_and_leporinB/leporin/...
$null = $ws.Cells.Replace(
    "/home/salamzade/zol_development/showcase_examples/Aflavus_aflatoxin_and_leporinB/leporinB_redo/",
    "/home/salamzade/zol_development/showcase_examples_redo/Aflavus_aflatoxin_and_leporinB/leporin/"
)

# 2. The data values (aggregate-bitscore ... copy-counts, columns C:AE) for the
#    two rows describing GCA_001695535.1_ASM169553v1's two gene clusters
#    (rows 197 and 198) were swapped between each other, while the
#    sample / gene-cluster-id columns (A:B) stayed put.
$row197 = $ws.Range("C197:AE197").Value()
$row198 = $ws.Range("C198:AE198").Value()
$ws.Range("C197:AE197").Value = $row198
$ws.Range("C198:AE198").Value = $row197
